$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from an existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$data = @{
    2  = @(8, 8)
    3  = @(2, 3)
    4  = @(6, 6)
    5  = @(7, 9)
    6  = @(8, 9)
    7  = @(7, 7)
    8  = @(8, 9)
    9  = @(4, 5)
    10 = @(8, 8)
    11 = @(4, 4)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
